$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect to allow writing, re-protect after.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure footnote (shared string).
$footnoteCell = $ws.Cells.Item(58, 1)
$footnoteCell.Value = $footnoteCell.Value2 -replace [regex]::Escape("2021-03-30"), "2021-03-31"

# Update Weight (D) and Percent Change (E) values for each holding row.
$ws.Cells.Item(2, 4).Value = 0.0163724787199547
$ws.Cells.Item(2, 5).Value = -0.005443863526114479
$ws.Cells.Item(3, 4).Value = 0.04992970140427561
$ws.Cells.Item(3, 5).Value = 0.01269601248981278
$ws.Cells.Item(4, 4).Value = 0.0147558904673457
$ws.Cells.Item(4, 5).Value = 0.0108755842027155
$ws.Cells.Item(5, 4).Value = 0.009724173326787362
$ws.Cells.Item(5, 5).Value = 0.005646679214842854
$ws.Cells.Item(6, 4).Value = 0.01583989131621713
$ws.Cells.Item(6, 5).Value = -0.007646976287357998
$ws.Cells.Item(7, 4).Value = 0.02080718809637247
$ws.Cells.Item(7, 5).Value = -0.000147655961609372
$ws.Cells.Item(8, 4).Value = 0.004344697365991678
$ws.Cells.Item(8, 5).Value = 0.006529752501316599
$ws.Cells.Item(9, 4).Value = 0.006670301792229336
$ws.Cells.Item(9, 5).Value = -0.01517022402540136
$ws.Cells.Item(10, 4).Value = 0.01413747461793677
$ws.Cells.Item(10, 5).Value = -0.008655666756829716
$ws.Cells.Item(11, 4).Value = 0.009041075644178274
$ws.Cells.Item(11, 5).Value = 0.00535030005061099
$ws.Cells.Item(12, 4).Value = 0.01495101454223451
$ws.Cells.Item(12, 5).Value = -0.0175760755508918
$ws.Cells.Item(13, 4).Value = 0.003048609394514961
$ws.Cells.Item(13, 5).Value = -0.004127579737335907
$ws.Cells.Item(14, 4).Value = 0.006219784162703802
$ws.Cells.Item(14, 5).Value = -0.02601156069364152
$ws.Cells.Item(15, 4).Value = 0.01464221516877563
$ws.Cells.Item(15, 5).Value = -0.01456499223200414
$ws.Cells.Item(16, 4).Value = 0.01073084359589549
$ws.Cells.Item(16, 5).Value = -0.002878289473684181
$ws.Cells.Item(17, 4).Value = 0.02157346681258156
$ws.Cells.Item(17, 5).Value = 0.02645214071448043
$ws.Cells.Item(18, 4).Value = 0.008690767464168504
$ws.Cells.Item(18, 5).Value = -0.002843152740483301
$ws.Cells.Item(19, 4).Value = 0.01725825050447091
$ws.Cells.Item(19, 5).Value = -0.003999757590449016
$ws.Cells.Item(20, 4).Value = 0.01212580096880247
$ws.Cells.Item(20, 5).Value = 0.007115902964959453
$ws.Cells.Item(21, 4).Value = 0.007432429627783342
$ws.Cells.Item(21, 5).Value = -0.01496908558411969
$ws.Cells.Item(22, 4).Value = 0.01392499528613573
$ws.Cells.Item(22, 5).Value = -0.01117245005257606
$ws.Cells.Item(23, 4).Value = 0.01983441123859762
$ws.Cells.Item(23, 5).Value = -0.01441537640149482
$ws.Cells.Item(24, 4).Value = 0.009955151858988244
$ws.Cells.Item(24, 5).Value = -0.009330667428353867
$ws.Cells.Item(25, 4).Value = 0.0212752570874081
$ws.Cells.Item(25, 5).Value = 0.003275283591628098
$ws.Cells.Item(26, 4).Value = 0.01157487778516329
$ws.Cells.Item(26, 5).Value = 0.008005218216318832
$ws.Cells.Item(27, 4).Value = 0.02058496890287681
$ws.Cells.Item(27, 5).Value = 0.03961156909319552
$ws.Cells.Item(28, 4).Value = 0.05486352966639307
$ws.Cells.Item(28, 5).Value = 0.01876563803169318
$ws.Cells.Item(29, 4).Value = 0.02151790384318274
$ws.Cells.Item(29, 5).Value = -0.002064220183486398
$ws.Cells.Item(30, 4).Value = 0.03095184236511035
$ws.Cells.Item(30, 5).Value = 0.02872228088701156
$ws.Cells.Item(31, 4).Value = 0.01573435435845901
$ws.Cells.Item(31, 5).Value = 0.01906079471216637
$ws.Cells.Item(32, 4).Value = 0.01369548753841894
$ws.Cells.Item(32, 5).Value = -0.009068649678062979
$ws.Cells.Item(33, 4).Value = 0.02036507028045605
$ws.Cells.Item(33, 5).Value = 0.01929743471227185
$ws.Cells.Item(34, 4).Value = 0.04013221764198144
$ws.Cells.Item(34, 5).Value = 0.007842773165499528
$ws.Cells.Item(35, 4).Value = 0.01140152132063896
$ws.Cells.Item(35, 5).Value = -0.007911936704506295
$ws.Cells.Item(36, 4).Value = 0.009772938003459745
$ws.Cells.Item(36, 5).Value = 0.01339076692574204
$ws.Cells.Item(37, 4).Value = 0.01198689354530541
$ws.Cells.Item(37, 5).Value = -0.01267893660531694
$ws.Cells.Item(38, 4).Value = 0.007447072104424914
$ws.Cells.Item(38, 5).Value = -0.02172481895984202
$ws.Cells.Item(39, 4).Value = 0.01185148332047051
$ws.Cells.Item(39, 5).Value = -0.007694280584765512
$ws.Cells.Item(40, 4).Value = 0.01793579189014011
$ws.Cells.Item(40, 5).Value = -0.001158972377824918
$ws.Cells.Item(41, 4).Value = 0.01734571315512459
$ws.Cells.Item(41, 5).Value = -0.009881139911212866
$ws.Cells.Item(42, 4).Value = 0.03273279895484094
$ws.Cells.Item(42, 5).Value = 0.02273612317646712
$ws.Cells.Item(43, 4).Value = 0.01133321155237805
$ws.Cells.Item(43, 5).Value = 0.0006575342465753309
$ws.Cells.Item(44, 4).Value = 0.02146985821670258
$ws.Cells.Item(44, 5).Value = -0.0008022652194431545
$ws.Cells.Item(45, 4).Value = 0.01341826100521844
$ws.Cells.Item(45, 5).Value = 0.01644159943879342
$ws.Cells.Item(46, 4).Value = 0.008410533993800547
$ws.Cells.Item(46, 5).Value = -0.006575265808617781
$ws.Cells.Item(47, 4).Value = 0.01371241790203575
$ws.Cells.Item(47, 5).Value = -0.005191350609233036
$ws.Cells.Item(48, 4).Value = 0.009840430669229492
$ws.Cells.Item(48, 5).Value = -0.01246524975338525
$ws.Cells.Item(49, 4).Value = 0.01469333310062254
$ws.Cells.Item(49, 5).Value = 0.006272855884472905
$ws.Cells.Item(50, 4).Value = 0.008345002374009579
$ws.Cells.Item(50, 5).Value = 0.00173114055529644
$ws.Cells.Item(51, 4).Value = 0.01108357039927909
$ws.Cells.Item(51, 5).Value = -0.001344688480501843
$ws.Cells.Item(52, 4).Value = 0.008815195831522227
$ws.Cells.Item(52, 5).Value = -0.008320078306619405
$ws.Cells.Item(53, 4).Value = 0.1420142432730485
$ws.Cells.Item(53, 5).Value = -0.00009850275807721243
$ws.Cells.Item(54, 4).Value = 0.04368360654135649
$ws.Cells.Item(54, 5).Value = 0.006359704909692176
$ws.Cells.Item(55, 5).Value = 0.003892382110985748

# Restore sheet protection.
$ws.Protect()

Write-Host "Edit complete."
